# Applies the commit "Added test methods in Login test and modified test
# methods in business" to the BusinessProfile worksheet.
#
# Net effect (per the OOXML diff):
#  - sheet2 (BusinessProfile) gets a new data row (row 3) that duplicates most
#    of row 2's values but uses the new "testAddToExistingAccount" /
#    "Add To Existing Account" test-case identifiers.
#  - two new columns (T: businessEmail, U: businessPassword) are appended,
#    with header cells on row 1, blank styled cells on row 2, and populated
#    hyperlinked cells on row 3.
#  - four new mailto hyperlinks are created (D3, E3, U3, T3 - in that order).
#  - the BusinessProfile sheet's selection moves to T5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "BusinessProfile"

function Copy-CellFormat($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy()
    $ws.Cells.Item($dstRow, $dstCol).PasteSpecial(-4122)  # xlPasteFormats
}

function Copy-CellText($srcRow, $srcCol, $dstRow, $dstCol) {
    # Force text format first so numeric-looking strings (e.g. "1",
    # "123456") are written back as shared strings instead of numbers -
    # the final Copy-CellFormat call (done separately) restores the real
    # number format/style afterwards.
    $ws.Cells.Item($dstRow, $dstCol).NumberFormat = "@"
    $ws.Cells.Item($dstRow, $dstCol).Value = $ws.Cells.Item($srcRow, $srcCol).Text
}

# ---------------------------------------------------------------------
# Row 3: new "Add To Existing Account" test case
# ---------------------------------------------------------------------

# A3 / B3 - brand new strings (must be written first so they land at the
# front of the newly appended shared-string block).
$ws.Cells.Item(3, 1).Value = "testAddToExistingAccount"
$ws.Cells.Item(3, 2).Value = "Add To Existing Account"

# C3 - reuse the "1" value/style (quote-prefixed text) from C2.
Copy-CellText  2, 3, 3, 3
Copy-CellFormat 2, 3, 3, 3

# D3 / E3 - same email/password text & style as D2/E2; hyperlinks added below.
Copy-CellText  2, 4, 3, 4
Copy-CellFormat 2, 4, 3, 4
Copy-CellText  2, 5, 3, 5
Copy-CellFormat 2, 5, 3, 5

# F3..M3 - straight carry-overs from row 2 (plain, unstyled cells).
Copy-CellText   2, 6, 3, 6
Copy-CellFormat 2, 6, 3, 6
Copy-CellText   2, 7, 3, 7
Copy-CellFormat 2, 7, 3, 7
Copy-CellText   2, 8, 3, 8
Copy-CellFormat 2, 8, 3, 8
Copy-CellText   2, 9, 3, 9
Copy-CellFormat 2, 9, 3, 9

Copy-CellText   2, 10, 3, 10
Copy-CellFormat 2, 10, 3, 10

Copy-CellText   2, 11, 3, 11
Copy-CellFormat 2, 11, 3, 11
Copy-CellText   2, 12, 3, 12
Copy-CellFormat 2, 12, 3, 12
Copy-CellText   2, 13, 3, 13
Copy-CellFormat 2, 13, 3, 13

# N3 - same "123456" value/style as J2/J3.
Copy-CellText   2, 10, 3, 14
Copy-CellFormat 2, 10, 3, 14

# P3 - carry-over from row 2 (column O/Q/R/S are intentionally left blank,
# matching the source diff).
Copy-CellText   2, 16, 3, 16
Copy-CellFormat 2, 16, 3, 16

# ---------------------------------------------------------------------
# New columns T (businessEmail) and U (businessPassword)
# ---------------------------------------------------------------------

# Header row (bold "Test Case Name" style, same as the rest of row 1).
$ws.Cells.Item(1, 20).Value = "businessEmail"
Copy-CellFormat 1, 1, 1, 20
$ws.Cells.Item(1, 21).Value = "businessPassword"
Copy-CellFormat 1, 1, 1, 21

# Row 2 placeholders - blank but styled like the Hyperlink cells.
Copy-CellFormat 2, 4, 2, 20
Copy-CellFormat 2, 4, 2, 21

# Row 3 values or the new columns.
$ws.Cells.Item(3, 20).Value = "shahul0100@gamil.com"
Copy-CellFormat 2, 4, 3, 20
$ws.Cells.Item(3, 21).Value = $ws.Cells.Item(2, 5).Text
Copy-CellFormat 2, 4, 3, 21

# ---------------------------------------------------------------------
# Hyperlinks - added in the same order as the target workbook
# (D3, E3, U3, T3) so the relationship ids line up.
# ---------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Cells.Item(3, 4), "mailto:vijayp@ideyalabs.com")
Copy-CellFormat 2, 4, 3, 4

$ws.Hyperlinks.Add($ws.Cells.Item(3, 5), "mailto:Admin@123")
Copy-CellFormat 2, 5, 3, 5

$ws.Hyperlinks.Add($ws.Cells.Item(3, 21), "mailto:Admin@123")
Copy-CellFormat 2, 4, 3, 21

$ws.Hyperlinks.Add($ws.Cells.Item(3, 20), "mailto:shahul0100@gamil.com")
Copy-CellFormat 2, 4, 3, 20

# ---------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------

$ws.Range("T5").Select()
